$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "d"
$ws.Range("C2").Value = "d"
$ws.Range("D2").Value = 10
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 3
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 269.1975
